$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 289, shifting the existing data (rows 289-309)
# down to rows 290-310.
$ws.Rows.Item(289).Insert()

# Populate the newly inserted row 289 with the new weekly price record.
$ws.Range("A289").Value = 4
$ws.Range("B289").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C289").Value = "Los Lagos"
$ws.Range("D289").Value = 44714
$ws.Range("E289").Value = 10
$ws.Range("F289").Value = 100112045
$ws.Range("G289").Value = "Zapallo"
$ws.Range("H289").Value = "Paine"
$ws.Range("I289").Value = "1a (cosecha)"
$ws.Range("J289").Value = 500
$ws.Range("K289").Value = 500
$ws.Range("L289").Value = 500
$ws.Range("M289").Value = 500
$ws.Range("N289").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O289").Value = "Región de O'Higgins"
$ws.Range("P289").Value = 500
$ws.Range("Q289").Value = 1
$ws.Range("R289").Value = "Hortaliza"
